$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("No Action")

$ws.Cells.Item(7,3).Value = 59.9
$ws.Cells.Item(7,4).Value = 58.33
$ws.Cells.Item(7,5).Value = 61.46
$ws.Cells.Item(7,6).Value = 63.54
$ws.Cells.Item(7,7).Value = 60.94
$ws.Cells.Item(7,8).Value = 61.98
$ws.Cells.Item(7,9).Value = 55.21
$ws.Cells.Item(7,10).Value = 57.81
$ws.Cells.Item(7,11).Value = 57.81
$ws.Cells.Item(7,12).Value = 65.62
$ws.Cells.Item(8,3).Value = 64.07
$ws.Cells.Item(8,4).Value = 58.01
$ws.Cells.Item(8,5).Value = 59.74
$ws.Cells.Item(8,6).Value = 65.37
$ws.Cells.Item(8,7).Value = 63.64
$ws.Cells.Item(8,8).Value = 60.17
$ws.Cells.Item(8,9).Value = 59.74
$ws.Cells.Item(8,10).Value = 61.04
$ws.Cells.Item(8,11).Value = 57.58
$ws.Cells.Item(8,12).Value = 62.77
$ws.Cells.Item(9,3).Value = 59.11
$ws.Cells.Item(9,4).Value = 61.71
$ws.Cells.Item(9,5).Value = 58.36
$ws.Cells.Item(9,6).Value = 63.94
$ws.Cells.Item(9,7).Value = 63.57
$ws.Cells.Item(9,8).Value = 59.48
$ws.Cells.Item(9,9).Value = 65.06
$ws.Cells.Item(9,10).Value = 56.51
$ws.Cells.Item(9,11).Value = 62.83
$ws.Cells.Item(9,12).Value = 58.36
$ws.Cells.Item(10,3).Value = 58.77
$ws.Cells.Item(10,4).Value = 60.71
$ws.Cells.Item(10,5).Value = 61.04
$ws.Cells.Item(10,6).Value = 61.69
$ws.Cells.Item(10,7).Value = 65.58
$ws.Cells.Item(10,8).Value = 62.01
$ws.Cells.Item(10,9).Value = 62.66
$ws.Cells.Item(10,10).Value = 64.29
$ws.Cells.Item(10,11).Value = 56.17
$ws.Cells.Item(10,12).Value = 62.01
$ws.Cells.Item(11,3).Value = 58.38
$ws.Cells.Item(11,4).Value = 63.01
$ws.Cells.Item(11,5).Value = 63.29
$ws.Cells.Item(11,6).Value = 63.58
$ws.Cells.Item(11,7).Value = 62.14
$ws.Cells.Item(11,8).Value = 58.38
$ws.Cells.Item(11,9).Value = 64.45
$ws.Cells.Item(11,10).Value = 58.67
$ws.Cells.Item(11,11).Value = 60.69
$ws.Cells.Item(11,12).Value = 61.27
$ws.Cells.Item(12,3).Value = 60.42
$ws.Cells.Item(12,4).Value = 65.89
$ws.Cells.Item(12,5).Value = 62.76
$ws.Cells.Item(12,6).Value = 60.68
$ws.Cells.Item(12,7).Value = 59.38
$ws.Cells.Item(12,8).Value = 60.68
$ws.Cells.Item(12,9).Value = 62.76
$ws.Cells.Item(12,10).Value = 64.06
$ws.Cells.Item(12,11).Value = 62.76
$ws.Cells.Item(12,12).Value = 63.28
$ws.Cells.Item(13,3).Value = 60.52
$ws.Cells.Item(13,4).Value = 64.54
$ws.Cells.Item(13,5).Value = 62.65
$ws.Cells.Item(13,6).Value = 63.83
$ws.Cells.Item(13,7).Value = 61.23
$ws.Cells.Item(13,8).Value = 64.78
$ws.Cells.Item(13,9).Value = 64.54
$ws.Cells.Item(13,10).Value = 60.52
$ws.Cells.Item(13,12).Value = 63.59
$ws.Cells.Item(14,3).Value = 64.21
$ws.Cells.Item(14,4).Value = 62.47
$ws.Cells.Item(14,5).Value = 63.12
$ws.Cells.Item(14,6).Value = 64.86
$ws.Cells.Item(14,7).Value = 63.99
$ws.Cells.Item(14,8).Value = 64.21
$ws.Cells.Item(14,9).Value = 63.56
$ws.Cells.Item(14,10).Value = 60.95
$ws.Cells.Item(14,11).Value = 62.26
$ws.Cells.Item(14,12).Value = 63.56
$ws.Cells.Item(15,3).Value = 60.2
$ws.Cells.Item(15,4).Value = 62
$ws.Cells.Item(15,5).Value = 63.6
$ws.Cells.Item(15,6).Value = 60
$ws.Cells.Item(15,7).Value = 62.6
$ws.Cells.Item(15,8).Value = 60.2
$ws.Cells.Item(15,9).Value = 62.6
$ws.Cells.Item(15,10).Value = 60.6
$ws.Cells.Item(15,11).Value = 64.8
$ws.Cells.Item(15,12).Value = 61.2
$ws.Cells.Item(16,3).Value = 66.36
$ws.Cells.Item(16,4).Value = 63.38
$ws.Cells.Item(16,5).Value = 62.45
$ws.Cells.Item(16,6).Value = 60.04
$ws.Cells.Item(16,7).Value = 62.64
$ws.Cells.Item(16,8).Value = 63.2
$ws.Cells.Item(16,9).Value = 60.78
$ws.Cells.Item(16,10).Value = 61.52
$ws.Cells.Item(16,11).Value = 67.66
$ws.Cells.Item(16,12).Value = 58.74
$ws.Cells.Item(17,3).Value = 64.76
$ws.Cells.Item(17,4).Value = 64.93
$ws.Cells.Item(17,5).Value = 63.89
$ws.Cells.Item(17,6).Value = 62.85
$ws.Cells.Item(17,7).Value = 60.07
$ws.Cells.Item(17,8).Value = 60.42
$ws.Cells.Item(17,9).Value = 63.02
$ws.Cells.Item(17,10).Value = 60.59
$ws.Cells.Item(17,11).Value = 61.46
$ws.Cells.Item(17,12).Value = 67.36
$ws.Cells.Item(18,3).Value = 60.16
$ws.Cells.Item(18,4).Value = 61.3
$ws.Cells.Item(18,5).Value = 63.25
$ws.Cells.Item(18,6).Value = 65.04
$ws.Cells.Item(18,7).Value = 64.39
$ws.Cells.Item(18,8).Value = 61.79
$ws.Cells.Item(18,9).Value = 65.2
$ws.Cells.Item(18,10).Value = 62.6
$ws.Cells.Item(18,11).Value = 60.16
$ws.Cells.Item(18,12).Value = 64.72
$ws.Cells.Item(19,3).Value = 62.02
$ws.Cells.Item(19,4).Value = 60.95
$ws.Cells.Item(19,5).Value = 60.64
$ws.Cells.Item(19,6).Value = 63.4
$ws.Cells.Item(19,7).Value = 63.55
$ws.Cells.Item(19,8).Value = 58.65
$ws.Cells.Item(19,9).Value = 60.64
$ws.Cells.Item(19,10).Value = 62.63
$ws.Cells.Item(19,11).Value = 59.88
$ws.Cells.Item(19,12).Value = 58.81
$ws.Cells.Item(20,3).Value = 59.25
$ws.Cells.Item(20,4).Value = 60.12
$ws.Cells.Item(20,5).Value = 55.78
$ws.Cells.Item(20,6).Value = 61.42
$ws.Cells.Item(20,7).Value = 60.69
$ws.Cells.Item(20,8).Value = 61.13
$ws.Cells.Item(20,9).Value = 56.21
$ws.Cells.Item(20,10).Value = 63.15
$ws.Cells.Item(20,11).Value = 64.74
$ws.Cells.Item(20,12).Value = 59.54
$ws.Cells.Item(21,3).Value = 63.42
$ws.Cells.Item(21,4).Value = 61.37
$ws.Cells.Item(21,5).Value = 58.08
$ws.Cells.Item(21,6).Value = 57.53
$ws.Cells.Item(21,7).Value = 55.07
$ws.Cells.Item(21,8).Value = 57.12
$ws.Cells.Item(21,9).Value = 61.37
$ws.Cells.Item(21,10).Value = 61.92
$ws.Cells.Item(21,11).Value = 58.9
$ws.Cells.Item(21,12).Value = 59.32

# K13 holds literal text "62,.88" (Excel kept it as text rather than a number
# because of the odd "62,.88" formatting). Enter it via a text formula then
# convert the formula to a plain value in place so the cell keeps its
# original numeric style (s="3") instead of acquiring a new text style.
$ws.Range("K13").Formula = "=""62,.88"""
$ws.Range("K13").Copy()
$ws.Range("K13").PasteSpecial(-4163)

# Recalculate the Min/Max/Mean formulas in columns M:O now that C:L are numbers.
$wb.Application.CalculateFull()

# Selection state on "No Action": the diff moves the live selection to S14.
$ws.Range("S14").Select()

# Activate the "No Action" tab (it becomes the selected tab; this also
# removes tabSelected from whichever sheet had it before, e.g. "Remove
# Incomplete Records").
$ws.Activate()
